$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new staff record (row 11): RAVI KUMAR C, Lab Instructor,
# photo path and unique id VEC-015-05-011.
# Write the unique id (column J) before the photo path (column C) so the
# shared-string table picks up the id ahead of the photo url, matching
# the authoring order used when this row was entered.
$ws.Range("A11").Value = "RAVI KUMAR C"
$ws.Range("B11").Value = "Lab Instructor"
$ws.Range("J11").Value = "VEC-015-05-011"
$ws.Range("C11").Value = "/static/images/profile_photos/015/VEC-015-05-011.webp"

# The Photo/unique_id cells in this sheet use the plain default cell
# style (no border/alignment override) rather than the column's general
# "blank" style, so normalize them explicitly.
$ws.Range("C11").Style = "Normal"
$ws.Range("J11").Style = "Normal"

# Leave the selection on the last-entered cell (Photo column), matching
# where editing finished for this row.
[void]$ws.Range("C11").Select()
